$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix mis-formatted numeric-looking text values (remove extra thousand-style comma) ---
$ws.Range("G2").Value = "'8235,50"
$ws.Range("F5").Value = "'4944,27"
$ws.Range("F6").Value = "'2242,03"

# --- Shift the GENERAL account codes for rows 6-8 down by one, inserting a new
#     row 9 for account 602612 with its own DEBIT amount ---
$ws.Range("B6").Value = "'601100"
$ws.Range("B7").Value = "'602201"
$ws.Range("B8").Value = "'602202"

# --- New row 9: new ledger line for account 602612 ---
$ws.Range("A9").Value = "13/07/2024"
$ws.Range("B9").Value = "'602612"
$ws.Range("D9").Value = "'878318"
$ws.Range("E9").Value = "Achat MB 878318"
$ws.Range("F9").Value = "'95,24"
